# Update the three headline COVID-19 stat cells on Sheet1 (row 2):
#   B2 = Total Deaths, C2 = Total Recovered, D2 = Total Infected
# Values keep the leading/trailing newline formatting used by the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "`n362,023`n"
$ws.Range("C2").Value = "`n2,579,678`n"
$ws.Range("D2").Value = "`n5,905,292 `n"

# Re-fit the row height back to its default so the multi-line text doesn't
# leave a stale auto-sized row height behind.
$ws.Rows(2).EntireRow.AutoFit()
